$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text columns (B: Coin name, C: Link, E: Volume change) - these remain text naturally
$textUpdates = @(
    @("E2", "  +5.31%  "),
    @("E3", "  +6.22%  "),
    @("E4", "  +0.48%  "),
    @("E5", "  +5.46%  "),
    @("E6", "  +7.95%  "),
    @("B7", "XRP"),
    @("C7", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"),
    @("E7", "  +23.28%  "),
    @("B8", "USDC"),
    @("C8", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"),
    @("E8", "  -0.46%  "),
    @("E9", "  +8.05%  "),
    @("E10", "  +14.63%  "),
    @("E11", "  +6.89%  "),
    @("E12", "  +6.00%  "),
    @("E13", "  +1.48%  "),
    @("E14", "  +6.41%  "),
    @("E15", "  +5.23%  "),
    @("E16", "  +8.48%  "),
    @("E17", "  +5.31%  "),
    @("E18", "  +7.21%  "),
    @("E19", "  +5.32%  "),
    @("E20", "  +6.28%  "),
    @("E21", "  +6.70%  "),
    @("E22", "  +6.92%  "),
    @("E23", "  +0.73%  "),
    @("E24", "  +4.98%  "),
    @("E25", "  +5.89%  "),
    @("E26", "  +6.86%  "),
    @("E27", "  -0.26%  "),
    @("E28", "  +4.97%  "),
    @("E29", "  +4.50%  "),
    @("E30", "  +7.66%  "),
    @("E31", "  -0.18%  "),
    @("E32", "  +8.46%  "),
    @("E33", "  +4.49%  "),
    @("E34", "  +5.71%  "),
    @("E35", "  +9.59%  "),
    @("E36", "  +10.37%  "),
    @("E37", "  +8.61%  "),
    @("E38", "  +2.76%  "),
    @("E39", "  +10.48%  "),
    @("E40", "  +7.10%  "),
    @("E41", "  +14.53%  "),
    @("E42", "  +3.97%  "),
    @("B43", "Stellar"),
    @("C43", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"),
    @("E43", "  +7.33%  "),
    @("B44", "Mantle"),
    @("C44", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"),
    @("E44", "  +7.26%  "),
    @("E45", "  +4.34%  "),
    @("E46", "  -0.62%  "),
    @("E47", "  +7.30%  "),
    @("E48", "  +13.03%  "),
    @("E49", "  +4.32%  "),
    @("E50", "  +0.65%  "),
    @("E51", "  +14.84%  ")
)
foreach ($item in $textUpdates) {
    $ws.Range($item[0]).Value = $item[1]
}

# Update price column (D) - force text format to preserve exact string representation
$priceUpdates = @(
    @("D2", "59.030.31"),
    @("D3", "2.535.82"),
    @("D5", "505.39"),
    @("D6", "159.94"),
    @("D7", "0.614"),
    @("D8", "0.995"),
    @("D9", "2.575.11"),
    @("D10", "6.24"),
    @("D14", "2.982.46"),
    @("D15", "58.972.94"),
    @("D16", "22.03"),
    @("D18", "2.570.88"),
    @("D19", "4.75"),
    @("D20", "333.91"),
    @("D21", "10.35"),
    @("D22", "6.07"),
    @("D24", "59.63"),
    @("D25", "0.417"),
    @("D26", "0.168"),
    @("D27", "0.997"),
    @("D28", "2.626.33"),
    @("D29", "7.57"),
    @("D30", "0.0₃0829"),
    @("D33", "154.75"),
    @("D35", "5.51"),
    @("D36", "3.95"),
    @("D37", "1.20"),
    @("D38", "0.860"),
    @("D39", "3.73"),
    @("D40", "1.44"),
    @("D41", "291.53"),
    @("D43", "0.101"),
    @("D44", "0.625"),
    @("D45", "0.0560"),
    @("D46", "0.991"),
    @("D47", "0.0239"),
    @("D48", "19.17"),
    @("D49", "4.81"),
    @("D51", "0.719")
)
foreach ($item in $priceUpdates) {
    $cell = $ws.Range($item[0])
    $cell.NumberFormat = "@"
    $cell.Value = $item[1]
    $cell.ClearFormats()
}

Write-Output "Done updating cryptos worksheet"
